$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Temáticas" table (column G) gets two new fields appended below the
# existing ID_Tematica / Nome_Tematica rows: a new "Quantidade_Tematica"
# attribute and a foreign key back to "ID_Candidato".
# Clone the formatting already used for the table (thin border, same as
# the rows above/alongside it) before writing the new text so the new
# cells match the rest of the sheet.
$ws.Range("G2").Copy() | Out-Null
$ws.Range("G4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("G3").Copy() | Out-Null
$ws.Range("G5").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("G4").Value = "Quantidade_Tematica"
$ws.Range("G5").Value = "ID_Candidato"

# Column G needs to be a bit wider now to comfortably fit the longer
# "Quantidade_Tematica" label.
$ws.Columns.Item(7).ColumnWidth = 20.6

# Leave the selection where the author ended up after the edit.
$ws.Range("G11").Select() | Out-Null
